# Generate Report for Handoff
#
# The "ht" (high/handoff-ready?) priority files (rows 4-7, the .xlf handoffs
# that were previously "low" priority) have had a fresh handoff xliff
# generated for both target languages. Reflect that in the per-language
# status tables:
#   - Priority goes from "low" to "ht" for rows 4-7
#   - Latest Handoff Datetime is bumped to the new generation timestamp

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

foreach ($row in 4..7) {
    $zhcn.Cells.Item($row, 5).Value = "ht"
    $dede.Cells.Item($row, 5).Value = "ht"

    $zhcn.Cells.Item($row, 8).Value = "2016-09-07 00:45:34"
    $dede.Cells.Item($row, 8).Value = "2016-09-07 00:45:40"
}
